$d = $word.ActiveDocument

# Locate the "PROMPT 5: object-position Matcher" section by finding its
# heading paragraph, then operate on the four body paragraphs that follow
# it (through the "Marcar completo: object-position ..." paragraph).
$headingIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "PROMPT 5: object-position Matcher") {
        $headingIndex = $i
        break
    }
}

# First body paragraph ("Implemente matchObjectPosition(value) em ...")
# only has its run text struck through; the paragraph mark itself keeps
# its original (empty) run properties.
$pImplemente = $d.Paragraphs($headingIndex + 1)
$rImplemente = $d.Range($pImplemente.Range.Start, $pImplemente.Range.End)
$rImplemente.Font.StrikeThrough = $true

# The remaining three paragraphs ("Mapear:", "Adicionar:",
# "Marcar completo:") have both their runs and their paragraph mark
# struck through.
$pMapear = $d.Paragraphs($headingIndex + 2)
$pMapear.Range.Font.StrikeThrough = $true

$pAdicionar = $d.Paragraphs($headingIndex + 3)
$pAdicionar.Range.Font.StrikeThrough = $true

$pMarcar = $d.Paragraphs($headingIndex + 4)
$pMarcar.Range.Font.StrikeThrough = $true
